# Update the "snapshot" sheet: remove the row for Порядин Павел (СПР)
# which has been replaced/returned, shifting subsequent rows up, then
# refresh the scraped_at (column K) timestamps for the remaining data rows.

$wb = $excel.ActiveWorkbook

$snapshot = $wb.Worksheets.Item("snapshot")

# Row 32 (team "СПР", player "Порядин Павел") is removed; everything below
# it shifts up by one row (native Excel row-delete semantics).
$snapshot.Rows.Item(32).Delete()

# New scrape timestamps for data rows 2..46 (post-shift).
$newTimestamps = @(
    "2025-11-14T07:02:32.702746+00:00",
    "2025-11-14T07:02:32.702785+00:00",
    "2025-11-14T07:02:32.702809+00:00",
    "2025-11-14T07:02:34.635737+00:00",
    "2025-11-14T07:02:34.635769+00:00",
    "2025-11-14T07:02:34.635791+00:00",
    "2025-11-14T07:02:37.060110+00:00",
    "2025-11-14T07:02:38.943707+00:00",
    "2025-11-14T07:02:41.006366+00:00",
    "2025-11-14T07:02:41.006405+00:00",
    "2025-11-14T07:02:46.108687+00:00",
    "2025-11-14T07:02:48.433474+00:00",
    "2025-11-14T07:02:50.302998+00:00",
    "2025-11-14T07:02:50.303028+00:00",
    "2025-11-14T07:02:50.303047+00:00",
    "2025-11-14T07:02:52.669790+00:00",
    "2025-11-14T07:02:54.591539+00:00",
    "2025-11-14T07:02:54.591571+00:00",
    "2025-11-14T07:02:56.966902+00:00",
    "2025-11-14T07:02:58.957548+00:00",
    "2025-11-14T07:02:58.957585+00:00",
    "2025-11-14T07:02:58.957611+00:00",
    "2025-11-14T07:02:58.957631+00:00",
    "2025-11-14T07:02:58.957654+00:00",
    "2025-11-14T07:03:01.371741+00:00",
    "2025-11-14T07:03:01.371772+00:00",
    "2025-11-14T07:03:03.728395+00:00",
    "2025-11-14T07:03:03.728427+00:00",
    "2025-11-14T07:03:03.728446+00:00",
    "2025-11-14T07:03:05.673334+00:00",
    "2025-11-14T07:03:05.673370+00:00",
    "2025-11-14T07:03:07.514020+00:00",
    "2025-11-14T07:03:07.514051+00:00",
    "2025-11-14T07:03:07.514070+00:00",
    "2025-11-14T07:03:07.514086+00:00",
    "2025-11-14T07:03:07.514102+00:00",
    "2025-11-14T07:03:07.514125+00:00",
    "2025-11-14T07:03:09.403645+00:00",
    "2025-11-14T07:03:09.403676+00:00",
    "2025-11-14T07:03:13.834739+00:00",
    "2025-11-14T07:03:13.834781+00:00",
    "2025-11-14T07:03:13.834806+00:00",
    "2025-11-14T07:03:13.834826+00:00",
    "2025-11-14T07:03:16.279323+00:00",
    "2025-11-14T07:03:16.279352+00:00"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = 2 + $i
    $snapshot.Cells.Item($row, 11).Value = $newTimestamps[$i]
}

# --- "returned" sheet -------------------------------------------------
# Порядин Павел (СПР) is now the only outstanding return; the four other
# previously-returned players (rows 3..6) are cleared out.
$returned = $wb.Worksheets.Item("returned")
$returned.Range("A3:G6").EntireRow.Delete()

$returned.Cells.Item(2, 1).Value = "СПР"
$returned.Cells.Item(2, 2).Value = "Спартак"
$returned.Cells.Item(2, 3).Value = "Порядин Павел"
$returned.Cells.Item(2, 4).Value = "1369_СПР_порядинпавел"
$returned.Cells.Item(2, 5).Value = "RETURN"
$returned.Cells.Item(2, 6).Value = "2025-11-14T15:03:16.786299+08:00"

# Column G holds a plain "yyyy-mm-dd" string; force text so Excel does not
# reinterpret it as a date serial number (it was stored as text before).
$changedDayCell = $returned.Cells.Item(2, 7)
$changedDayCell.NumberFormat = "@"
$changedDayCell.Value = "2025-11-14"

